$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of the "Betarraga" block
# (row 215), pushing every subsequent record down by one row
# (old 215..253 -> new 216..254).
$ws.Rows(215).Insert()

$ws.Range("A215").Value = 10
$ws.Range("B215").Value = "Vega Modelo de Temuco"
$ws.Range("C215").Value = "La Araucanía"
$ws.Range("D215").Value = 44476
$ws.Range("E215").Value = 9
$ws.Range("F215").Value = 100114014
$ws.Range("G215").Value = "Betarraga"
$ws.Range("H215").Value = "Sin especificar"
$ws.Range("I215").Value = "Primera"
$ws.Range("J215").Value = 40
$ws.Range("K215").Value = 8000
$ws.Range("L215").Value = 8000
$ws.Range("M215").Value = 8000
$ws.Range("N215").Value = "$/docena de paquetes"
$ws.Range("O215").Value = "Provincia de Cautín"
$ws.Range("P215").Value = 667
$ws.Range("Q215").Value = 12
$ws.Range("R215").Value = "Hortaliza"
